$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style/formatting of the existing header cell (H1) onto the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I (I0) and J (IF), rows 2-35
$data = @(
    @(7, 8),
    @(7, 7),
    @(6, 7),
    @(9, 10),
    @(6, 6),
    @(9, 9),
    @(8, 8),
    @(9, 9),
    @(5, 7),
    @(7, 7),
    @(8, 9),
    @(6, 6),
    @(7, 7),
    @(7, 9),
    @(5, 6),
    @(7, 7),
    @(7, 7),
    @(4, 5),
    @(12, 12),
    @(7, 8),
    @(8, 8),
    @(7, 7),
    @(5, 6),
    @(5, 5),
    @(8, 8),
    @(7, 8),
    @(8, 8),
    @(8, 9),
    @(7, 7),
    @(8, 8),
    @(7, 7),
    @(4, 4),
    @(6, 6),
    @(3, 3)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
